# Auto-generated edit script: updates market-data derived columns (H:N)
# on several sheets to reflect refreshed pricing data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 333.66666
$ws.Range("I5").Value = 194
$ws.Range("J5").Value = 508.25
$ws.Range("K5").Value = 194
$ws.Range("L5").Value = 508.25
$ws.Range("M5").Value = -79
$ws.Range("N5").Value = -738.25
$ws.Range("H15").Value = 931.3111
$ws.Range("I15").Value = 931.3111
$ws.Range("K15").Value = 2793.9333
$ws.Range("M15").Value = -2624.9333
$ws.Range("H19").Value = 1316546.1
$ws.Range("I19").Value = 2393041
$ws.Range("K19").Value = 2393041
$ws.Range("M19").Value = -2392866
$ws.Range("H41").Value = 741.1177
$ws.Range("I41").Value = 757
$ws.Range("J41").Value = 730
$ws.Range("K41").Value = 757
$ws.Range("L41").Value = 730
$ws.Range("M41").Value = -317
$ws.Range("N41").Value = -1610
$ws.Range("H92").Value = 1268.9524
$ws.Range("I92").Value = 1087.8823
$ws.Range("J92").Value = 2038.5
$ws.Range("K92").Value = 1087.8823
$ws.Range("L92").Value = 2038.5
$ws.Range("M92").Value = 160.1177
$ws.Range("N92").Value = -4534.5
$ws.Range("H96").Value = 627.46155
$ws.Range("I96").Value = 246.44444
$ws.Range("J96").Value = 1484.75
$ws.Range("K96").Value = 739.33332
$ws.Range("L96").Value = 4454.25
$ws.Range("M96").Value = 633.66668
$ws.Range("N96").Value = -7200.25
$ws.Range("H100").Value = 50000996
$ws.Range("I100").Value = 50000996
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 50000996
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -50000455
$ws.Range("N100").Value = ""
$ws.Range("H107").Value = 1424.2174
$ws.Range("I107").Value = 1707
$ws.Range("J107").Value = 984.3333
$ws.Range("K107").Value = 1707
$ws.Range("L107").Value = 984.3333
$ws.Range("M107").Value = 213
$ws.Range("N107").Value = -4824.3333
$ws.Range("H112").Value = 1622.766
$ws.Range("J112").Value = 1659.591
$ws.Range("L112").Value = 4978.772999999999
$ws.Range("N112").Value = -7194.772999999999
$ws.Range("H116").Value = 406680.8
$ws.Range("I116").Value = 1251713.5
$ws.Range("J116").Value = 9018.352999999999
$ws.Range("K116").Value = 1251713.5
$ws.Range("L116").Value = 9018.352999999999
$ws.Range("M116").Value = -1248271.5
$ws.Range("N116").Value = -15902.353
$ws.Range("H132").Value = 12787631
$ws.Range("I132").Value = 16669547
$ws.Range("J132").Value = 528949.1
$ws.Range("K132").Value = 50008641
$ws.Range("L132").Value = 1586847.3
$ws.Range("M132").Value = -50006111
$ws.Range("N132").Value = -1591907.3
$ws.Range("H135").Value = 685.04
$ws.Range("I135").Value = 286.3889
$ws.Range("J135").Value = 1710.1428
$ws.Range("K135").Value = 2577.5001
$ws.Range("L135").Value = 15391.2852
$ws.Range("M135").Value = -42.50009999999975
$ws.Range("N135").Value = -20461.2852
$ws.Range("H138").Value = 3536.74
$ws.Range("I138").Value = 911.25806
$ws.Range("J138").Value = 4716.304
$ws.Range("K138").Value = 2733.77418
$ws.Range("L138").Value = 14148.912
$ws.Range("M138").Value = 2406.22582
$ws.Range("N138").Value = -24428.912
$ws.Range("H141").Value = 5812.6665
$ws.Range("I141").Value = 6507.4326
$ws.Range("J141").Value = 2599.375
$ws.Range("K141").Value = 19522.2978
$ws.Range("L141").Value = 7798.125
$ws.Range("M141").Value = -14342.2978
$ws.Range("N141").Value = -18158.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 914.881
$ws.Range("I61").Value = 799.88
$ws.Range("J61").Value = 1084
$ws.Range("K61").Value = 799.88
$ws.Range("L61").Value = 1084
$ws.Range("M61").Value = -587.88
$ws.Range("N61").Value = -1508
$ws.Range("H74").Value = 2078.4443
$ws.Range("I74").Value = 2014.4117
$ws.Range("J74").Value = 2350.5833
$ws.Range("K74").Value = 2014.4117
$ws.Range("L74").Value = 2350.5833
$ws.Range("M74").Value = -1140.4117
$ws.Range("N74").Value = -4098.5833
$ws.Range("H77").Value = 2078.4443
$ws.Range("I77").Value = 2014.4117
$ws.Range("J77").Value = 2350.5833
$ws.Range("K77").Value = 10072.0585
$ws.Range("L77").Value = 11752.9165
$ws.Range("M77").Value = -5704.058500000001
$ws.Range("N77").Value = -20488.9165
$ws.Range("H97").Value = 527.36365
$ws.Range("I97").Value = 552.76666
$ws.Range("J97").Value = 273.33334
$ws.Range("K97").Value = 552.76666
$ws.Range("L97").Value = 273.33334
$ws.Range("M97").Value = -56.76666
$ws.Range("N97").Value = -1265.33334
$ws.Range("H102").Value = 1854.25
$ws.Range("I102").Value = 1772.3334
$ws.Range("J102").Value = 2100
$ws.Range("K102").Value = 1772.3334
$ws.Range("L102").Value = 2100
$ws.Range("M102").Value = -150.3334
$ws.Range("N102").Value = -5344
$ws.Range("H132").Value = 1938.2115
$ws.Range("I132").Value = 1280.5
$ws.Range("J132").Value = 3723.4285
$ws.Range("K132").Value = 3841.5
$ws.Range("L132").Value = 11170.2855
$ws.Range("M132").Value = -1311.5
$ws.Range("N132").Value = -16230.2855
$ws.Range("H136").Value = 914.881
$ws.Range("I136").Value = 799.88
$ws.Range("J136").Value = 1084
$ws.Range("K136").Value = 2399.64
$ws.Range("L136").Value = 3252
$ws.Range("M136").Value = 150.3600000000001
$ws.Range("N136").Value = -8352

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2475.6223
$ws.Range("I134").Value = 1394.4231
$ws.Range("J134").Value = 3955.158
$ws.Range("K134").Value = 4183.2693
$ws.Range("L134").Value = 11865.474
$ws.Range("M134").Value = -1648.2693
$ws.Range("N134").Value = -16935.474
$ws.Range("H140").Value = 48905.832
$ws.Range("J140").Value = 48905.832
$ws.Range("L140").Value = 48905.832
$ws.Range("N140").Value = -59265.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8335389.5
$ws.Range("I31").Value = 1178.6586
$ws.Range("J31").Value = 26319740
$ws.Range("K31").Value = 1178.6586
$ws.Range("L31").Value = 26319740
$ws.Range("M31").Value = -883.6586
$ws.Range("N31").Value = -26320330
$ws.Range("H34").Value = 8335389.5
$ws.Range("I34").Value = 1178.6586
$ws.Range("J34").Value = 26319740
$ws.Range("K34").Value = 1178.6586
$ws.Range("L34").Value = 26319740
$ws.Range("M34").Value = -976.6586
$ws.Range("N34").Value = -26320144
$ws.Range("H58").Value = 1656.253
$ws.Range("I58").Value = 1487.0143
$ws.Range("J58").Value = 2567.5386
$ws.Range("K58").Value = 1487.0143
$ws.Range("L58").Value = 2567.5386
$ws.Range("M58").Value = -1284.0143
$ws.Range("N58").Value = -2973.5386
$ws.Range("H107").Value = 1053.091
$ws.Range("I107").Value = 630.2
$ws.Range("K107").Value = 630.2
$ws.Range("M107").Value = 1289.8
$ws.Range("H132").Value = 2704.1428
$ws.Range("I132").Value = 2354.3333
$ws.Range("K132").Value = 7062.999899999999
$ws.Range("M132").Value = -4532.999899999999
$ws.Range("H134").Value = 3375.3462
$ws.Range("I134").Value = 3690.3235
$ws.Range("K134").Value = 11070.9705
$ws.Range("M134").Value = -8535.970499999999
$ws.Range("H136").Value = 1656.253
$ws.Range("I136").Value = 1487.0143
$ws.Range("J136").Value = 2567.5386
$ws.Range("K136").Value = 4461.0429
$ws.Range("L136").Value = 7702.6158
$ws.Range("M136").Value = -1911.0429
$ws.Range("N136").Value = -12802.6158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 275.75
$ws.Range("I7").Value = 151.5
$ws.Range("K7").Value = 454.5
$ws.Range("M7").Value = -342.5
$ws.Range("H131").Value = 11628716
$ws.Range("J131").Value = 991.5484
$ws.Range("L131").Value = 2974.6452
$ws.Range("N131").Value = -13054.6452

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 26886.727
$ws.Range("J43").Value = 26886.727
$ws.Range("L43").Value = 26886.727
$ws.Range("N43").Value = -27188.727
$ws.Range("H57").Value = 38000
$ws.Range("J57").Value = 37500
$ws.Range("L57").Value = 37500
$ws.Range("N57").Value = -39140
$ws.Range("H80").Value = 125001600
$ws.Range("I80").Value = 250000000
$ws.Range("J80").Value = 3200
$ws.Range("K80").Value = 250000000
$ws.Range("L80").Value = 3200
$ws.Range("M80").Value = -249999002
$ws.Range("N80").Value = -5196
$ws.Range("H83").Value = 125001600
$ws.Range("I83").Value = 250000000
$ws.Range("J83").Value = 3200
$ws.Range("K83").Value = 1250000000
$ws.Range("L83").Value = 16000
$ws.Range("M83").Value = -1249995008
$ws.Range("N83").Value = -25984
$ws.Range("H107").Value = 6944917
$ws.Range("I107").Value = 277.3
$ws.Range("J107").Value = 18519316
$ws.Range("K107").Value = 277.3
$ws.Range("L107").Value = 18519316
$ws.Range("M107").Value = 1642.7
$ws.Range("N107").Value = -18523156
$ws.Range("H126").Value = 3140.31
$ws.Range("I126").Value = 2909.5518
$ws.Range("J126").Value = 4684.615
$ws.Range("K126").Value = 8728.6554
$ws.Range("L126").Value = 14053.845
$ws.Range("M126").Value = -6258.6554
$ws.Range("N126").Value = -18993.845
$ws.Range("H132").Value = 2077.0852
$ws.Range("I132").Value = 1293.4688
$ws.Range("J132").Value = 3748.8
$ws.Range("K132").Value = 3880.4064
$ws.Range("L132").Value = 11246.4
$ws.Range("M132").Value = -1350.4064
$ws.Range("N132").Value = -16306.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3906.4375
$ws.Range("J2").Value = 3906.4375
$ws.Range("L2").Value = 3906.4375
$ws.Range("N2").Value = -4130.4375
$ws.Range("H45").Value = 39996
$ws.Range("J45").Value = 39996
$ws.Range("L45").Value = 39996
$ws.Range("N45").Value = -40810
$ws.Range("H112").Value = 31000
$ws.Range("J112").Value = 31000
$ws.Range("L112").Value = 31000
$ws.Range("N112").Value = -33954
$ws.Range("H132").Value = 22021.176
$ws.Range("I132").Value = 42668.5
$ws.Range("J132").Value = 10759
$ws.Range("K132").Value = 128005.5
$ws.Range("L132").Value = 32277
$ws.Range("M132").Value = -125475.5
$ws.Range("N132").Value = -37337
$ws.Range("H136").Value = 2281.889
$ws.Range("I136").Value = 1024.1892
$ws.Range("J136").Value = 8098.75
$ws.Range("K136").Value = 3072.5676
$ws.Range("L136").Value = 24296.25
$ws.Range("M136").Value = -522.5676000000003
$ws.Range("N136").Value = -29396.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""
$ws.Range("H132").Value = 9011169
$ws.Range("I132").Value = 1617.64
$ws.Range("K132").Value = 4852.92
$ws.Range("M132").Value = -2322.92
$ws.Range("H136").Value = 1930.7091
$ws.Range("I136").Value = 645.75
$ws.Range("J136").Value = 5357.2666
$ws.Range("K136").Value = 1937.25
$ws.Range("L136").Value = 16071.7998
$ws.Range("M136").Value = 612.75
$ws.Range("N136").Value = -21171.7998
